# Add a new row (row 6) for user Hitesh to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting from the row above so the new row matches existing styling exactly
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data row: EIN=2121, Name=Hitesh, Tower=Compute
$ws.Range("A6").Value = 2121
$ws.Range("B6").Value = "Hitesh"
$ws.Range("C6").Value = "Compute"

# Update the selected/active cell as in the diff
$ws.Range("D13").Select()

$wb.Save()
